# Adds a new "forecast column" (AH, made on 2020-05-13) and a new row 46
# (observation date 2020-05-27) to both the "cases" and "deaths" sheets of
# the forecasts table, mirroring the staircase pattern already present in
# the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data: the new AH-column forecast values (rows 33-45) per sheet,
# plus the new B32 value and the new AH46 / row-46 value.
# ---------------------------------------------------------------------

$sheetNames = @("cases", "deaths")

$casesAH = @{
    33 = 52323
    34 = 53383
    35 = 54465
    36 = 55705
    37 = 56930
    38 = 58097
    39 = 59021
    40 = 60428
    41 = 61299
    42 = 62148
    43 = 63048
    44 = 63962
    45 = 64830
}
$casesB32 = 51097
$casesAH46 = 65805

$deathsAH = @{
    33 = 4235
    34 = 4329
    35 = 4418
    36 = 4510
    37 = 4589
    38 = 4671
    39 = 4753
    40 = 4833
    41 = 4917
    42 = 5002
    43 = 5079
    44 = 5154
    45 = 5224
}
$deathsB32 = 4118
$deathsAH46 = 5295

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "cases") {
        $ahValues = $casesAH
        $b32 = $casesB32
        $ah46 = $casesAH46
    } else {
        $ahValues = $deathsAH
        $b32 = $deathsB32
        $ah46 = $deathsAH46
    }

    # Column AH is column 34 (A=1 ... AG=33, AH=34).
    $col = 34

    # --- Row 1: header for the new forecast-date column. It reuses the
    # same text label that the existing header cells use (a literal date
    # string), so force text entry and strip the format change back off so
    # the cell ends up with the default style, same as its neighbours.
    $ws.Cells.Item(1, $col).NumberFormat = "@"
    $ws.Cells.Item(1, $col).Value = "2020-05-13"
    $ws.Cells.Item(1, $col).ClearFormats()

    # --- Rows 2-32: blank placeholder cells in the new column (matches the
    # staircase layout where the forecast hadn't been made yet for those
    # observation dates). Materialize them as empty cells (no value) rather
    # than leaving them absent.
    for ($r = 2; $r -le 32; $r++) {
        $ws.Cells.Item($r, $col).NumberFormat = "General"
        $ws.Cells.Item($r, $col).ClearFormats()
    }

    # --- Rows 33-45: actual forecast numbers for the new column.
    foreach ($r in 33..45) {
        $ws.Cells.Item($r, $col).Value = $ahValues[$r]
    }

    # --- B32 gains a forecast value that was previously blank.
    $ws.Cells.Item(32, 2).Value = $b32

    # --- New row 46 (observation date 2020-05-27).
    $ws.Cells.Item(46, 1).NumberFormat = "@"
    $ws.Cells.Item(46, 1).Value = "2020-05-27"
    $ws.Cells.Item(46, 1).ClearFormats()

    for ($c = 2; $c -le 33; $c++) {
        $ws.Cells.Item(46, $c).NumberFormat = "General"
        $ws.Cells.Item(46, $c).ClearFormats()
    }

    $ws.Cells.Item(46, $col).Value = $ah46
}
